# Regenerate the lattice-multiplication practice problems.
# Every cell in the 5x3 table gets a brand new "A x B" problem plus the
# matching lattice scaffold lines (top multiplier digits, the divider, and
# the two left-edge multiplicand digits).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# [char]11 is the vertical-tab character Word uses for a manual line break
# (<w:br/>) inside a single paragraph/run.
$nl = [char]11

# row, col, top line, 2nd line, 3rd line, 4th line, 5th line
$newProblems = @(
    @(1, 1, "89 x 69", "  6    9", "  ----", "8|    |", "9|    |"),
    @(1, 2, "12 x 52", "  5    2", "  ----", "1|    |", "2|    |"),
    @(1, 3, "68 x 13", "  1    3", "  ----", "6|    |", "8|    |"),

    @(2, 1, "77 x 53", "  5    3", "  ----", "7|    |", "7|    |"),
    @(2, 2, "25 x 29", "  2    9", "  ----", "2|    |", "5|    |"),
    @(2, 3, "79 x 70", "  7    0", "  ----", "7|    |", "9|    |"),

    @(3, 1, "60 x 15", "  1    5", "  ----", "6|    |", "0|    |"),
    @(3, 2, "15 x 28", "  2    8", "  ----", "1|    |", "5|    |"),
    @(3, 3, "77 x 68", "  6    8", "  ----", "7|    |", "7|    |"),

    @(4, 1, "53 x 73", "  7    3", "  ----", "5|    |", "3|    |"),
    @(4, 2, "82 x 23", "  2    3", "  ----", "8|    |", "2|    |"),
    @(4, 3, "53 x 29", "  2    9", "  ----", "5|    |", "3|    |"),

    @(5, 1, "67 x 74", "  7    4", "  ----", "6|    |", "7|    |"),
    @(5, 2, "38 x 25", "  2    5", "  ----", "3|    |", "8|    |"),
    @(5, 3, "76 x 11", "  1    1", "  ----", "7|    |", "6|    |")
)

foreach ($p in $newProblems) {
    $row = $p[0]
    $col = $p[1]
    $cellLines = $p[2..6]
    $newText = [string]::Join($nl, $cellLines)

    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}

Write-Output "Updated $($newProblems.Count) lattice-multiplication cells"
